$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header "Категория" -> "Категории"
$ws.Range("C1").Value = "Категории"

# Populate the "Этап" (stage) and "Результат" (result) columns with real values
# instead of the placeholder text that used to be duplicated across E:H.
$ws.Range("G2").Value = "школьный"
$ws.Range("H2").Value = "призер"

$ws.Range("G3").Value = "отборочный"
$ws.Range("H3").Value = "победитель"

$ws.Range("G4").Value = "школьный"
$ws.Range("H4").Value = "победитель"

$ws.Range("G5").Value = "школьный"
$ws.Range("H5").Value = "победитель"

# The "Баллы" (points) column is no longer part of the export - drop it entirely.
$ws.Columns.Item(9).Delete()

# Resize the stage/result columns to fit their new contents.
$ws.Columns.Item(7).EntireColumn.AutoFit()
$ws.Columns.Item(8).EntireColumn.AutoFit()
